# edit.ps1
# Commit: "case with 380 kV done"
#
# The slack-bus voltage setpoint (column B, "vm_pu" of the slack bus) was changed
# from 1.05 pu to 1.02 pu for this network case. Because the bus-voltage results
# (res_bus/vm_pu) are the output of a power-flow solve, changing the slack setpoint
# required the whole case to be re-solved, which updated the voltage-magnitude
# results for every other bus as well (columns C-F and I-N; column G is the slack
# bus itself and stays at 1 pu; column H has no bus mapped to it). This updates
# those recomputed per-unit voltage magnitudes for data rows 2-25 (bus scenarios
# 0-23) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)


# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.011882166886027
$ws.Cells.Item(2, 4).Value = 1.043197011191348
$ws.Cells.Item(2, 5).Value = 1.013998921470454
$ws.Cells.Item(2, 6).Value = 1.046742216697039
$ws.Cells.Item(2, 9).Value = 1.036288903430969
$ws.Cells.Item(2, 10).Value = 1.017128081965357
$ws.Cells.Item(2, 11).Value = 1.04597155682114
$ws.Cells.Item(2, 12).Value = 1.016858390300359
$ws.Cells.Item(2, 13).Value = 1.049506801078452
$ws.Cells.Item(2, 14).Value = 1.009826501662289

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.012779032203262
$ws.Cells.Item(3, 4).Value = 1.043710875753346
$ws.Cells.Item(3, 5).Value = 1.014758044684493
$ws.Cells.Item(3, 6).Value = 1.047437083147122
$ws.Cells.Item(3, 9).Value = 1.036360220882024
$ws.Cells.Item(3, 10).Value = 1.017658679692866
$ws.Cells.Item(3, 11).Value = 1.046297184953277
$ws.Cells.Item(3, 12).Value = 1.017422491881995
$ws.Cells.Item(3, 13).Value = 1.050013674919614
$ws.Cells.Item(3, 14).Value = 1.010005446472735

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.013360223167837
$ws.Cells.Item(4, 4).Value = 1.044043380790664
$ws.Cells.Item(4, 5).Value = 1.015250375530645
$ws.Cells.Item(4, 6).Value = 1.047887108883981
$ws.Cells.Item(4, 9).Value = 1.036405042810243
$ws.Cells.Item(4, 10).Value = 1.018002278928982
$ws.Cells.Item(4, 11).Value = 1.046507173346474
$ws.Cells.Item(4, 12).Value = 1.017787971845068
$ws.Cells.Item(4, 13).Value = 1.050341379377689
$ws.Cells.Item(4, 14).Value = 1.010121228696372

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.013604760228603
$ws.Cells.Item(5, 4).Value = 1.044183163784159
$ws.Cells.Item(5, 5).Value = 1.015457619682565
$ws.Cells.Item(5, 6).Value = 1.048076393245307
$ws.Cells.Item(5, 9).Value = 1.036423567873966
$ws.Cells.Item(5, 10).Value = 1.01814679084343
$ws.Cells.Item(5, 11).Value = 1.046595279720683
$ws.Cells.Item(5, 12).Value = 1.017941730548223
$ws.Cells.Item(5, 13).Value = 1.050479078094131
$ws.Cells.Item(5, 14).Value = 1.010169901398223

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.013645831024387
$ws.Cells.Item(6, 4).Value = 1.044206633775232
$ws.Cells.Item(6, 5).Value = 1.015492432571367
$ws.Cells.Item(6, 6).Value = 1.048108180338409
$ws.Cells.Item(6, 9).Value = 1.036426659637983
$ws.Cells.Item(6, 10).Value = 1.018171058679259
$ws.Cells.Item(6, 11).Value = 1.046610062991294
$ws.Cells.Item(6, 12).Value = 1.017967553779879
$ws.Cells.Item(6, 13).Value = 1.05050219425569
$ws.Cells.Item(6, 14).Value = 1.010178073619854

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.013363489882555
$ws.Cells.Item(7, 4).Value = 1.044045248588995
$ws.Cells.Item(7, 5).Value = 1.015253143685098
$ws.Cells.Item(7, 6).Value = 1.047889637746339
$ws.Cells.Item(7, 9).Value = 1.036405291593889
$ws.Cells.Item(7, 10).Value = 1.018004209658929
$ws.Cells.Item(7, 11).Value = 1.046508351308991
$ws.Cells.Item(7, 12).Value = 1.017790025942072
$ws.Cells.Item(7, 13).Value = 1.050343219585296
$ws.Cells.Item(7, 14).Value = 1.010121879072402

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.012185088047307
$ws.Cells.Item(8, 4).Value = 1.043370672744206
$ws.Cells.Item(8, 5).Value = 1.014255236306557
$ws.Cells.Item(8, 6).Value = 1.046976965673494
$ws.Cells.Item(8, 9).Value = 1.036313279407703
$ws.Cells.Item(8, 10).Value = 1.017307343981573
$ws.Cells.Item(8, 11).Value = 1.046081751397282
$ws.Cells.Item(8, 12).Value = 1.017048933316463
$ws.Cells.Item(8, 13).Value = 1.04967815759951
$ws.Cells.Item(8, 14).Value = 1.009886978078182

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.010115230981826
$ws.Cells.Item(9, 4).Value = 1.042182083169535
$ws.Cells.Item(9, 5).Value = 1.012505502363995
$ws.Cells.Item(9, 6).Value = 1.045371894766127
$ws.Cells.Item(9, 9).Value = 1.036141028339403
$ws.Cells.Item(9, 10).Value = 1.016081472781086
$ws.Cells.Item(9, 11).Value = 1.045324625207025
$ws.Cells.Item(9, 12).Value = 1.015746676489813
$ws.Cells.Item(9, 13).Value = 1.048504200279087
$ws.Cells.Item(9, 14).Value = 1.009473019019609

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.008739866177731
$ws.Cells.Item(10, 4).Value = 1.04138988789033
$ws.Cells.Item(10, 5).Value = 1.011344963952945
$ws.Cells.Item(10, 6).Value = 1.044304123581206
$ws.Cells.Item(10, 9).Value = 1.036019442739917
$ws.Cells.Item(10, 10).Value = 1.015265701362877
$ws.Cells.Item(10, 11).Value = 1.044816338919899
$ws.Cells.Item(10, 12).Value = 1.014881028459589
$ws.Cells.Item(10, 13).Value = 1.047720312529139
$ws.Cells.Item(10, 14).Value = 1.009197051256456

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.008145410967451
$ws.Cells.Item(11, 4).Value = 1.041046933269902
$ws.Cells.Item(11, 5).Value = 1.010843870002949
$ws.Cells.Item(11, 6).Value = 1.043842336833083
$ws.Cells.Item(11, 9).Value = 1.035965203681352
$ws.Cells.Item(11, 10).Value = 1.014912827605719
$ws.Cells.Item(11, 11).Value = 1.044595427290222
$ws.Cells.Item(11, 12).Value = 1.014506807321925
$ws.Cells.Item(11, 13).Value = 1.047380608957135
$ws.Cells.Item(11, 14).Value = 1.009077561126036

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.007924768397675
$ws.Cells.Item(12, 4).Value = 1.040919557604022
$ws.Cells.Item(12, 5).Value = 1.010657957429504
$ws.Cells.Item(12, 6).Value = 1.043670896039033
$ws.Cells.Item(12, 9).Value = 1.035944818557976
$ws.Cells.Item(12, 10).Value = 1.014781809900538
$ws.Cells.Item(12, 11).Value = 1.044513249288648
$ws.Cells.Item(12, 12).Value = 1.014367897853135
$ws.Cells.Item(12, 13).Value = 1.047254388620293
$ws.Cells.Item(12, 14).Value = 1.009033178496756

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.007972089510418
$ws.Cells.Item(13, 4).Value = 1.04094687949454
$ws.Cells.Item(13, 5).Value = 1.010697826505036
$ws.Cells.Item(13, 6).Value = 1.04370766667864
$ws.Cells.Item(13, 9).Value = 1.035949202007696
$ws.Cells.Item(13, 10).Value = 1.014809911121447
$ws.Cells.Item(13, 11).Value = 1.044530882224437
$ws.Cells.Item(13, 12).Value = 1.014397690180608
$ws.Cells.Item(13, 13).Value = 1.047281465057403
$ws.Cells.Item(13, 14).Value = 1.009042698653964

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.008127169219331
$ws.Cells.Item(14, 4).Value = 1.041036404080416
$ws.Cells.Item(14, 5).Value = 1.010828497995046
$ws.Cells.Item(14, 6).Value = 1.043828163687214
$ws.Cells.Item(14, 9).Value = 1.035963523497703
$ws.Cells.Item(14, 10).Value = 1.014901996506251
$ws.Cells.Item(14, 11).Value = 1.044588636901739
$ws.Cells.Item(14, 12).Value = 1.014495323119035
$ws.Cells.Item(14, 13).Value = 1.047370176334914
$ws.Cells.Item(14, 14).Value = 1.009073892414331

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.008222740790929
$ws.Cells.Item(15, 4).Value = 1.04109156491589
$ws.Cells.Item(15, 5).Value = 1.010909037688698
$ws.Cells.Item(15, 6).Value = 1.043902417520923
$ws.Cells.Item(15, 9).Value = 1.035972315884637
$ws.Cells.Item(15, 10).Value = 1.014958740715473
$ws.Cells.Item(15, 11).Value = 1.044624205402668
$ws.Cells.Item(15, 12).Value = 1.014555490345495
$ws.Cells.Item(15, 13).Value = 1.047424829132613
$ws.Cells.Item(15, 14).Value = 1.00909311210914

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.008779341219539
$ws.Cells.Item(16, 4).Value = 1.041412650351327
$ws.Cells.Item(16, 5).Value = 1.011378250112031
$ws.Cells.Item(16, 6).Value = 1.044334782964421
$ws.Cells.Item(16, 9).Value = 1.036023008941801
$ws.Cells.Item(16, 10).Value = 1.015289128127424
$ws.Cells.Item(16, 11).Value = 1.04483098295568
$ws.Cells.Item(16, 12).Value = 1.014905877246543
$ws.Cells.Item(16, 13).Value = 1.047742851914224
$ws.Cells.Item(16, 14).Value = 1.00920498158328

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.009128773839035
$ws.Cells.Item(17, 4).Value = 1.041614079411703
$ws.Cells.Item(17, 5).Value = 1.011672957849221
$ws.Cells.Item(17, 6).Value = 1.044606147439556
$ws.Cells.Item(17, 9).Value = 1.036054381789619
$ws.Cells.Item(17, 10).Value = 1.015496468780235
$ws.Cells.Item(17, 11).Value = 1.04496047060172
$ws.Cells.Item(17, 12).Value = 1.015125829936746
$ws.Cells.Item(17, 13).Value = 1.047942266998372
$ws.Cells.Item(17, 14).Value = 1.009275156211256

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.009332696750791
$ws.Cells.Item(18, 4).Value = 1.041731576289905
$ws.Cells.Item(18, 5).Value = 1.011844993444063
$ws.Cells.Item(18, 6).Value = 1.044764484118643
$ws.Cells.Item(18, 9).Value = 1.03607252741144
$ws.Cells.Item(18, 10).Value = 1.015617441796319
$ws.Cells.Item(18, 11).Value = 1.045035919295762
$ws.Cells.Item(18, 12).Value = 1.015254183408325
$ws.Cells.Item(18, 13).Value = 1.048058555774667
$ws.Cells.Item(18, 14).Value = 1.009316088406426

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.009402246926962
$ws.Cells.Item(19, 4).Value = 1.041771640759152
$ws.Cells.Item(19, 5).Value = 1.011903676417667
$ws.Cells.Item(19, 6).Value = 1.044818481982297
$ws.Cells.Item(19, 9).Value = 1.036078688516612
$ws.Cells.Item(19, 10).Value = 1.015658696323355
$ws.Cells.Item(19, 11).Value = 1.045061631856462
$ws.Cells.Item(19, 12).Value = 1.015297958558671
$ws.Cells.Item(19, 13).Value = 1.048098202657519
$ws.Cells.Item(19, 14).Value = 1.009330045301601

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.009091272179166
$ws.Cells.Item(20, 4).Value = 1.04159246726894
$ws.Cells.Item(20, 5).Value = 1.011641324257851
$ws.Cells.Item(20, 6).Value = 1.04457702695794
$ws.Cells.Item(20, 9).Value = 1.036051031664722
$ws.Cells.Item(20, 10).Value = 1.015474219490858
$ws.Cells.Item(20, 11).Value = 1.044946585988944
$ws.Cells.Item(20, 12).Value = 1.015102225004869
$ws.Cells.Item(20, 13).Value = 1.047920874391011
$ws.Cells.Item(20, 14).Value = 1.009267627080463

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.008081497557247
$ws.Cells.Item(21, 4).Value = 1.041010040945632
$ws.Cells.Item(21, 5).Value = 1.010790012521267
$ws.Cells.Item(21, 6).Value = 1.043792677890168
$ws.Cells.Item(21, 9).Value = 1.035959312748193
$ws.Cells.Item(21, 10).Value = 1.01487487811936
$ws.Cells.Item(21, 11).Value = 1.04457163293187
$ws.Cells.Item(21, 12).Value = 1.014466570076887
$ws.Cells.Item(21, 13).Value = 1.04734405414364
$ws.Cells.Item(21, 14).Value = 1.009064706589069

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.007447564926404
$ws.Cells.Item(22, 4).Value = 1.040643922437217
$ws.Cells.Item(22, 5).Value = 1.010256009755318
$ws.Cells.Item(22, 6).Value = 1.043300033175481
$ws.Cells.Item(22, 9).Value = 1.035900266968022
$ws.Cells.Item(22, 10).Value = 1.014498369038922
$ws.Cells.Item(22, 11).Value = 1.044335182418614
$ws.Cells.Item(22, 12).Value = 1.014067446844401
$ws.Cells.Item(22, 13).Value = 1.046981157945773
$ws.Cells.Item(22, 14).Value = 1.008937130320726

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.007783533869373
$ws.Cells.Item(23, 4).Value = 1.040838000826362
$ws.Cells.Item(23, 5).Value = 1.010538975616547
$ws.Cells.Item(23, 6).Value = 1.043561144617569
$ws.Cells.Item(23, 9).Value = 1.035931698644426
$ws.Cells.Item(23, 10).Value = 1.014697932854643
$ws.Cells.Item(23, 11).Value = 1.044460595360376
$ws.Cells.Item(23, 12).Value = 1.014278978183303
$ws.Cells.Item(23, 13).Value = 1.047173556877813
$ws.Cells.Item(23, 14).Value = 1.00900476001069

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.009108217241411
$ws.Cells.Item(24, 4).Value = 1.041602232845515
$ws.Cells.Item(24, 5).Value = 1.011655617691317
$ws.Cells.Item(24, 6).Value = 1.044590185082418
$ws.Cells.Item(24, 9).Value = 1.036052545916749
$ws.Cells.Item(24, 10).Value = 1.015484272879868
$ws.Cells.Item(24, 11).Value = 1.044952860092908
$ws.Cells.Item(24, 12).Value = 1.015112890876025
$ws.Cells.Item(24, 13).Value = 1.047930540872096
$ws.Cells.Item(24, 14).Value = 1.009271029169025

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.010649544134275
$ws.Cells.Item(25, 4).Value = 1.042489336770077
$ws.Cells.Item(25, 5).Value = 1.012956808764195
$ws.Cells.Item(25, 6).Value = 1.045786453277501
$ws.Cells.Item(25, 9).Value = 1.036186752917476
$ws.Cells.Item(25, 10).Value = 1.016398134645756
$ws.Cells.Item(25, 11).Value = 1.01608290177176
$ws.Cells.Item(25, 12).Value = 1.009580038498504
$ws.Cells.Item(25, 13).Value = 1.0488079252043314
$ws.Cells.Item(25, 14).Value = 1.0095800384985503

